$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-03-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-24 Friday", 2)

# Update each arithmetic cell in the table by explicit (row, col) position
# to avoid ambiguity from duplicate expressions appearing more than once.
$t = $d.Tables.Item(1)

$cellValues = @(
    @(1,1,"56-24="), @(1,2,"18+76="), @(1,3,"24+12="), @(1,4,"87-80="), @(1,5,"80-61="),
    @(2,1,"77-7="), @(2,2,"39+1="), @(2,3,"14+52="), @(2,4,"54-12="), @(2,5,"51-43="),
    @(3,1,"58+36="), @(3,2,"33+34="), @(3,3,"17+4="), @(3,4,"18+24="), @(3,5,"53-46="),
    @(4,1,"40+59="), @(4,2,"33+54="), @(4,3,"63+35="), @(4,4,"35+44="), @(4,5,"14+0="),
    @(5,1,"71-44="), @(5,2,"82-53="), @(5,3,"2+24="), @(5,4,"42+48="), @(5,5,"39-18="),
    @(6,1,"92-12="), @(6,2,"78-55="), @(6,3,"82+10="), @(6,4,"20+5="), @(6,5,"17+19="),
    @(7,1,"91-0="), @(7,2,"79-35="), @(7,3,"50-33="), @(7,4,"31-4="), @(7,5,"28+71="),
    @(8,1,"88-38="), @(8,2,"87+10="), @(8,3,"62-50="), @(8,4,"6+20="), @(8,5,"68+20="),
    @(9,1,"96-5="), @(9,2,"79-31="), @(9,3,"26-24="), @(9,4,"94-36="), @(9,5,"4+1="),
    @(10,1,"41+38="), @(10,2,"16+39="), @(10,3,"1+38="), @(10,4,"86-43="), @(10,5,"67+13="),
    @(11,1,"65-27="), @(11,2,"62+35="), @(11,3,"25+69="), @(11,4,"10+77="), @(11,5,"31+17="),
    @(12,1,"61-58="), @(12,2,"2+6="), @(12,3,"47-15="), @(12,4,"3+88="), @(12,5,"35+60="),
    @(13,1,"55+22="), @(13,2,"79-1="), @(13,3,"21+58="), @(13,4,"64-31="), @(13,5,"32+20="),
    @(14,1,"52-13="), @(14,2,"8+31="), @(14,3,"49-49="), @(14,4,"52+24="), @(14,5,"28-4="),
    @(15,1,"22+1="), @(15,2,"38-23="), @(15,3,"94-38="), @(15,4,"41+7="), @(15,5,"86-66="),
    @(16,1,"84+4="), @(16,2,"68+2="), @(16,3,"64-54="), @(16,4,"67-67="), @(16,5,"64-4="),
    @(17,1,"31+14="), @(17,2,"49+34="), @(17,3,"74-63="), @(17,4,"61+21="), @(17,5,"90+3="),
    @(18,1,"62+30="), @(18,2,"22-10="), @(18,3,"52-38="), @(18,4,"10+54="), @(18,5,"70-9="),
    @(19,1,"75-62="), @(19,2,"18+71="), @(19,3,"8+45="), @(19,4,"1+32="), @(19,5,"59+13="),
    @(20,1,"64-3="), @(20,2,"92-44="), @(20,3,"94-21="), @(20,4,"37-36="), @(20,5,"51+22=")
)

foreach ($entry in $cellValues) {
    $row = $entry[0]
    $col = $entry[1]
    $newText = $entry[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
